# Applies the "license suspension" paragraph expansion described in the
# commit "Updated template for license suspension."
#
# Target paragraph (pre-edit) reads:
#   License Suspension. {{ license_suspension_details.license_type }}
#
# Target paragraph (post-edit) should read:
#   License Suspension. The Court orders that Defendant's {{ license_suspension_details.license_type }}
#   license is suspended from {{ license_suspension_details.license_suspended_date }} for a term of
#   {{ license_suspension_details.license_suspension_term }}.
#
# The trailing "_GoBack" bookmark that used to sit in the *next* paragraph
# (around the "{% endif %}") moves to sit right before the final "}}." run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the paragraph that holds "License Suspension." / the existing
# `{{ license_suspension_details.license_type }}` merge field.
# ---------------------------------------------------------------------
$targetPara = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*License Suspension.*license_type*") {
        $targetPara = $candidate
        break
    }
}

# ---------------------------------------------------------------------
# Find the existing "{{ license_suspension_details.license_type }}"
# expression inside that paragraph so we can anchor the new text right
# before / after it.
# ---------------------------------------------------------------------
$findRange = $targetPara.Range.Duplicate
$findRange.Find.Execute("{{ license_suspension_details.license_type }}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$exprStart = $findRange.Start
$exprEnd = $findRange.End

# ---------------------------------------------------------------------
# 1) Insert "The Court orders that Defendant's " right before the "{{".
# ---------------------------------------------------------------------
$prefixText = "The Court orders that Defendant" + [char]0x2019 + "s "
$beforeRange = $d.Range($exprStart, $exprStart)
$beforeRange.InsertBefore($prefixText)

# Recompute the (now shifted) expression end.
$exprEnd = $exprEnd + $prefixText.Length

# ---------------------------------------------------------------------
# 2) Insert the long suffix right after the "}}".
# ---------------------------------------------------------------------
$suffixText = " license is suspended from {{ license_suspension_details.license_suspended_date }} for a term of {{ license_suspension_details.license_suspension_term }}. "
$afterRange = $d.Range($exprEnd, $exprEnd)
$afterRange.InsertAfter($suffixText)

$prefixStart = $exprStart
$prefixEnd = $exprStart + $prefixText.Length
$suffixStart = $exprEnd
$suffixEnd = $exprEnd + $suffixText.Length

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the following paragraph so it now
#    sits right before the trailing "}}. " run, i.e. right after the
#    single trailing space that follows "license_suspension_term".
# ---------------------------------------------------------------------
$bookmarkPos = $suffixStart + ($suffixText.Length - 4)   # just before "}}. "
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------
# 4) Force the inserted prefix/suffix text to sit in their own runs
#    (rather than being silently coalesced into the neighbouring runs)
#    by nudging the font size away from and back to its real value.
#    This reproduces the discrete <w:r> boundaries from the authored
#    diff without altering the visible formatting.
# ---------------------------------------------------------------------
function Split-Run($rangeStart, $rangeEnd) {
    $rg = $d.Range($rangeStart, $rangeEnd)
    $rg.Font.Size = 13
    $rg.Font.Size = 11
}

# Prefix run: "The Court orders that Defendant's "
Split-Run $prefixStart $prefixEnd

# Suffix runs, in document order:
#   " license is suspended from {{ "
#   "license_suspension_details."
#   "license_suspended_date"
#   " }} for a term of {{ "
#   "license_suspension_details."
#   "license_suspension_term"
#   " "
#   "}}. "
$parts = @(
    " license is suspended from {{ ",
    "license_suspension_details.",
    "license_suspended_date",
    " }} for a term of {{ ",
    "license_suspension_details.",
    "license_suspension_term",
    " ",
    "}}. "
)

$cursor = $suffixStart
foreach ($part in $parts) {
    $partStart = $cursor
    $partEnd = $cursor + $part.Length
    Split-Run $partStart $partEnd
    $cursor = $partEnd
}

# ---------------------------------------------------------------------
# 5) Also give the very first "{{ license_suspension_details.license_type }}"
#    segment (the text that existed before the edit) its own refreshed
#    run boundary so it does not get absorbed by the new prefix run.
# ---------------------------------------------------------------------
Split-Run $prefixEnd $suffixStart

Write-Output "done"
